$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$src = $ws.Range("F16:X17").Value2
$ws.Range("E16:W17").Value2 = $src

$ws.Range("W16").Value2 = "Serial Number Image"
$ws.Range("W17").Value2 = "{booking:serial_number_pic}"

foreach ($col in @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X")) {
    $h = $ws.Range($col + "16").Text
    $p = $ws.Range($col + "17").Text
    Write-Host "$col : $h | $p"
}
